$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = "Hi, I am Clever Brain Technologies Assistant. How can I help you today ?"
